$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "fdec4c2032023fdccca6"
$ws.Range("B4").Value = "b40965018b5ec67fbf48"
$ws.Range("B5").Value = "0129c39d3efddff7cd09"
$ws.Range("B6").Value = "db42be9ae6c37108dc0f"
